# The deck ships two DrawingML themes:
#   ppt/theme/theme1.xml -> clrScheme "Office"   (used by the Notes Master)
#   ppt/theme/theme2.xml -> clrScheme "Integral" (used by the Slide Master /
#                                                  all slides, and referenced
#                                                  as the presentation theme)
#
# The target revision swaps the two palettes: the theme driving the slides
# (theme2.xml) becomes the plain "Office" palette, while the Notes Master
# palette becomes "Integral". The relationship targets (file names) are left
# untouched - only the RGB values inside the theme parts change.
#
# PowerPoint exposes the live, editable color palette for the theme that
# backs the slides through Slide.ThemeColorScheme (12 slots, in
# dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink order). Writing to it rewrites
# ppt/theme/theme2.xml in place.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)
$colors = $slide.ThemeColorScheme

# Target palette = the presentation's original "Office Theme" colors.
# The .RGB property takes a COLORREF (0xBBGGRR) - the bytes of the usual
# RRGGBB hex code in reverse order.
$colors.Item(1).RGB  = 0x000000   # dk1      000000
$colors.Item(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$colors.Item(3).RGB  = 0x6A5444   # dk2      44546A
$colors.Item(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$colors.Item(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$colors.Item(6).RGB  = 0x317DED   # accent2  ED7D31
$colors.Item(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$colors.Item(8).RGB  = 0x00C0FF   # accent4  FFC000
$colors.Item(9).RGB  = 0xC47244   # accent5  4472C4
$colors.Item(10).RGB = 0x47AD70   # accent6  70AD47
$colors.Item(11).RGB = 0xC16305   # hlink    0563C1
$colors.Item(12).RGB = 0x724F95   # folHlink 954F72
